$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2025")

# Insert a new row at row 51, shifting existing rows 51+ down by one.
$ws.Rows.Item(51).Insert()

# Populate new row 51 with the new benchmark entry.
$ws.Cells.Item(51, 1).Value = "Eigenbau 6 AI CUDA"
$ws.Cells.Item(51, 2).Value = "i3-6100"
$ws.Cells.Item(51, 3).Value = 3700
$ws.Cells.Item(51, 4).Value = 68.150000000000006
$ws.Cells.Item(51, 5).Value = 48.103000000000002
$ws.Cells.Item(51, 6).Value = 118.447
$ws.Cells.Item(51, 7).Value = 223.80500000000001
$ws.Cells.Item(51, 8).Value = 216.554
$ws.Cells.Item(51, 10).Value = 45773
$ws.Cells.Item(51, 12).Value = "x86-64"
$ws.Cells.Item(51, 13).Formula = "=D51*1000"
$ws.Cells.Item(51, 14).Formula = "=E51*1000"
$ws.Cells.Item(51, 15).Formula = "=F51*1000"
$ws.Cells.Item(51, 16).Formula = "=G51*1000"
$ws.Cells.Item(51, 17).Formula = "=H51*1000"
